$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Formed the consolidated report: fill in the "Absent" (column H) values
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
